# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for the data rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K values keyed by row number (row 1 is header; data starts row 2)
$kValues = @{
    2  = 0
    3  = 0
    4  = 3
    5  = 1
    6  = 2
    7  = 2
    8  = 1
    9  = 0
    10 = 0
    11 = 3
    12 = 0
    13 = 1
    14 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
